$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "110÷6=" "421÷3="
Replace-Text "649÷5=" "267÷5="
Replace-Text "540÷5=" "643÷9="
Replace-Text "194÷7=" "330÷9="
Replace-Text "563÷5=" "190÷3="
Replace-Text "223÷7=" "324÷3="
Replace-Text "454÷2=" "607÷4="
Replace-Text "234÷9=" "869÷2="
Replace-Text "407÷9=" "648÷2="
Replace-Text "548÷7=" "623÷3="
Replace-Text "577÷9=" "946÷9="
Replace-Text "643÷2=" "845÷4="
Replace-Text "271÷6=" "711÷2="
Replace-Text "652÷6=" "365÷7="
Replace-Text "245÷5=" "437÷3="
Replace-Text "127÷7=" "409÷7="
Replace-Text "880÷2=" "320÷3="
Replace-Text "410÷8=" "359÷7="
Replace-Text "903÷6=" "903÷3="
Replace-Text "586÷2=" "479÷2="
Replace-Text "803÷2=" "446÷2="
Replace-Text "452÷2=" "862÷7="
Replace-Text "509÷9=" "436÷5="
Replace-Text "864÷8=" "432÷2="
Replace-Text "652÷4=" "579÷8="
